# results_BRS.xlsx - "ISIP 2020 Overview" leaderboard refresh
#
# The sheet holds six parallel mini-tables (rank | col1 | col2 | score),
# one per "Subject" group, in the column bands B:D, G:I, L:N, Q:S, V:X and
# AA:AC, for data rows 11-20 and 30-40. This refresh pushes an updated
# results snapshot into those bands: some previously-empty placeholder
# cells (-1) receive real figures, some figures are revised, and a few
# previously-populated cells fall back to the -1 placeholder again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("L11").Value = 364
$ws.Range("M11").Value = 358
$ws.Range("N11").Value = 461.55
$ws.Range("V11").Value = 380
$ws.Range("W11").Value = 387
$ws.Range("X11").Value = 504.86

# Row 12
$ws.Range("L12").Value = 414
$ws.Range("M12").Value = 427
$ws.Range("N12").Value = 421.18
$ws.Range("V12").Value = 434
$ws.Range("W12").Value = 445
$ws.Range("X12").Value = 442.04
$ws.Range("AA12").Value = 500
$ws.Range("AB12").Value = 366
$ws.Range("AC12").Value = 424.18

# Row 13
$ws.Range("B13").Value = 811
$ws.Range("C13").Value = 412
$ws.Range("D13").Value = 502.02
$ws.Range("L13").Value = 504
$ws.Range("M13").Value = 424
$ws.Range("N13").Value = 373.57
$ws.Range("Q13").Value = 459
$ws.Range("R13").Value = 421
$ws.Range("S13").Value = 364.21
$ws.Range("V13").Value = 515
$ws.Range("W13").Value = 426
$ws.Range("X13").Value = 364.49
$ws.Range("AA13").Value = 604
$ws.Range("AB13").Value = 410
$ws.Range("AC13").Value = 356.82

# Row 14
$ws.Range("B14").Value = 735
$ws.Range("C14").Value = 483
$ws.Range("D14").Value = 430.01
$ws.Range("G14").Value = -1
$ws.Range("H14").Value = -1
$ws.Range("I14").Value = -1
$ws.Range("L14").Value = 557
$ws.Range("M14").Value = 358
$ws.Range("N14").Value = 315.02
$ws.Range("Q14").Value = 453
$ws.Range("R14").Value = 316
$ws.Range("S14").Value = 287.27999999999997
$ws.Range("V14").Value = 534
$ws.Range("W14").Value = 346
$ws.Range("X14").Value = 297.45999999999998
$ws.Range("AA14").Value = 682
$ws.Range("AB14").Value = 347
$ws.Range("AC14").Value = 311.2

# Row 15
$ws.Range("B15").Value = 643
$ws.Range("C15").Value = 454
$ws.Range("D15").Value = 362.89
$ws.Range("G15").Value = -1
$ws.Range("H15").Value = -1
$ws.Range("I15").Value = -1
$ws.Range("L15").Value = 527
$ws.Range("M15").Value = 279
$ws.Range("N15").Value = 241.5
$ws.Range("Q15").Value = 359
$ws.Range("R15").Value = 252
$ws.Range("S15").Value = 219.81
$ws.Range("V15").Value = 478
$ws.Range("W15").Value = 280
$ws.Range("X15").Value = 244.37
$ws.Range("AA15").Value = 666
$ws.Range("AB15").Value = 240
$ws.Range("AC15").Value = 263.97000000000003

# Row 16
$ws.Range("G16").Value = 392
$ws.Range("H16").Value = 270
$ws.Range("I16").Value = 230.28
$ws.Range("Q16").Value = 250
$ws.Range("R16").Value = 269
$ws.Range("S16").Value = 179.59
$ws.Range("V16").Value = 397
$ws.Range("W16").Value = 258
$ws.Range("X16").Value = 208.78
$ws.Range("AA16").Value = 576
$ws.Range("AB16").Value = 180
$ws.Range("AC16").Value = 215.04

# Row 17
$ws.Range("G17").Value = 272
$ws.Range("H17").Value = 344
$ws.Range("I17").Value = 172.05
$ws.Range("Q17").Value = 185
$ws.Range("R17").Value = 349
$ws.Range("S17").Value = 149.61000000000001
$ws.Range("V17").Value = 324
$ws.Range("W17").Value = 318
$ws.Range("X17").Value = 171.97
$ws.Range("AA17").Value = 474
$ws.Range("AB17").Value = 229
$ws.Range("AC17").Value = 158.82

# Row 18
$ws.Range("B18").Value = 857
$ws.Range("C18").Value = 315
$ws.Range("D18").Value = 180.2
$ws.Range("G18").Value = 270
$ws.Range("H18").Value = 481
$ws.Range("I18").Value = 115.37
$ws.Range("Q18").Value = 213
$ws.Range("R18").Value = 458
$ws.Range("S18").Value = 117.43
$ws.Range("V18").Value = 308
$ws.Range("W18").Value = 402
$ws.Range("X18").Value = 139.97999999999999

# Row 19
$ws.Range("B19").Value = 831
$ws.Range("C19").Value = 507
$ws.Range("D19").Value = 105.26
$ws.Range("Q19").Value = 294
$ws.Range("R19").Value = 508
$ws.Range("S19").Value = 89.02
$ws.Range("V19").Value = 344
$ws.Range("W19").Value = 464
$ws.Range("X19").Value = 112.66

# Row 20
$ws.Range("B20").Value = 628
$ws.Range("C20").Value = 526
$ws.Range("D20").Value = 22.19
$ws.Range("V20").Value = 508
$ws.Range("W20").Value = 516
$ws.Range("X20").Value = 37.28

# Row 30
$ws.Range("Q30").Value = 362
$ws.Range("R30").Value = 317
$ws.Range("S30").Value = 545.41

# Row 31
$ws.Range("Q31").Value = 359
$ws.Range("R31").Value = 408
$ws.Range("S31").Value = 479.97
$ws.Range("AA31").Value = 360
$ws.Range("AB31").Value = 334
$ws.Range("AC31").Value = 563.13

# Row 32
$ws.Range("Q32").Value = 446
$ws.Range("R32").Value = 440
$ws.Range("S32").Value = 406.48

# Row 33
$ws.Range("L33").Value = 570
$ws.Range("M33").Value = 432
$ws.Range("N33").Value = 364.98
$ws.Range("Q33").Value = 520
$ws.Range("R33").Value = 382
$ws.Range("S33").Value = 336.72
$ws.Range("V33").Value = 409
$ws.Range("W33").Value = 394
$ws.Range("X33").Value = 465.2

# Row 34
$ws.Range("L34").Value = 612
$ws.Range("M34").Value = 322
$ws.Range("N34").Value = 301.68
$ws.Range("Q34").Value = 519
$ws.Range("R34").Value = 292
$ws.Range("S34").Value = 287.22000000000003
$ws.Range("V34").Value = 512
$ws.Range("W34").Value = 460
$ws.Range("X34").Value = 381.46

# Row 35
$ws.Range("B35").Value = 546
$ws.Range("C35").Value = 216
$ws.Range("D35").Value = 251.18
$ws.Range("G35").Value = 542
$ws.Range("H35").Value = 240
$ws.Range("I35").Value = 261.01
$ws.Range("L35").Value = 534
$ws.Range("M35").Value = 236
$ws.Range("N35").Value = 251.13
$ws.Range("Q35").Value = 449
$ws.Range("R35").Value = 236
$ws.Range("S35").Value = 246.12

# Row 36
$ws.Range("B36").Value = 450
$ws.Range("C36").Value = 171
$ws.Range("D36").Value = 209.57
$ws.Range("G36").Value = 448
$ws.Range("H36").Value = 179
$ws.Range("I36").Value = 220.79
$ws.Range("L36").Value = 414
$ws.Range("M36").Value = 228
$ws.Range("N36").Value = 206.87
$ws.Range("Q36").Value = 344
$ws.Range("R36").Value = 258
$ws.Range("S36").Value = 200.32

# Row 37
$ws.Range("B37").Value = 357
$ws.Range("C37").Value = 209
$ws.Range("D37").Value = 175.5
$ws.Range("G37").Value = 335
$ws.Range("H37").Value = 202
$ws.Range("I37").Value = 186.02
$ws.Range("L37").Value = 325
$ws.Range("M37").Value = 306
$ws.Range("N37").Value = 168.35
$ws.Range("Q37").Value = 290
$ws.Range("R37").Value = 327
$ws.Range("S37").Value = 165

# Row 38
$ws.Range("B38").Value = 319
$ws.Range("C38").Value = 312
$ws.Range("D38").Value = 139.13
$ws.Range("G38").Value = 266
$ws.Range("H38").Value = 302
$ws.Range("I38").Value = 151.85
$ws.Range("L38").Value = 311
$ws.Range("M38").Value = 426
$ws.Range("N38").Value = 131.35
$ws.Range("Q38").Value = 319
$ws.Range("R38").Value = 434
$ws.Range("S38").Value = 118.28

# Row 39
$ws.Range("B39").Value = 348
$ws.Range("C39").Value = 414
$ws.Range("D39").Value = 103.11
$ws.Range("G39").Value = 270
$ws.Range("H39").Value = 426
$ws.Range("I39").Value = 117.95
$ws.Range("L39").Value = 379
$ws.Range("M39").Value = 500
$ws.Range("N39").Value = 100.2
$ws.Range("Q39").Value = 436
$ws.Range("R39").Value = 506
$ws.Range("S39").Value = 56.32

# Row 40
$ws.Range("B40").Value = 532
$ws.Range("C40").Value = 480
$ws.Range("D40").Value = 30.36
$ws.Range("G40").Value = 355
$ws.Range("H40").Value = 510
$ws.Range("I40").Value = 84.52

# Reproduce the author's last on-screen selection at save time
$ws.Range("AJ30").Select()
